# Fruta / hortaliza, semanal
# Insert 3 new daily-report rows (date 44585) above the current top row (875)
# for this product block. All rows from 875 downward shift down by 3, which
# automatically reproduces the duplicated tail rows (923-925) expected by the
# diff, since old rows 920-922 land exactly on 923-925 after the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 875 (existing rows 875-922 shift to 878-925)
$ws.Rows("875:877").Insert()

# Row 875 (new): Hass / Especial
$ws.Range("A875").Value = 10
$ws.Range("B875").Value = "Vega Modelo de Temuco"
$ws.Range("C875").Value = "La Araucanía"
$ws.Range("D875").Value = 44585
$ws.Range("E875").Value = 9
$ws.Range("F875").Value = "Fruta"
$ws.Range("G875").Value = 100106
$ws.Range("H875").Value = "Oleaginosos"
$ws.Range("I875").Value = 100106002
$ws.Range("J875").Value = "Palta"
$ws.Range("K875").Value = "Hass"
$ws.Range("L875").Value = "Especial"
$ws.Range("M875").Value = 140
$ws.Range("N875").Value = 3500
$ws.Range("O875").Value = 3500
$ws.Range("P875").Value = 3500
$ws.Range("Q875").Value = "`$/kilo (en bandeja de 18 kilos)"
$ws.Range("R875").Value = "Provincia de Quillota"
$ws.Range("S875").Value = 3500
$ws.Range("T875").Value = 1

# Row 876 (new): Hass / Primera
$ws.Range("A876").Value = 10
$ws.Range("B876").Value = "Vega Modelo de Temuco"
$ws.Range("C876").Value = "La Araucanía"
$ws.Range("D876").Value = 44585
$ws.Range("E876").Value = 9
$ws.Range("F876").Value = "Fruta"
$ws.Range("G876").Value = 100106
$ws.Range("H876").Value = "Oleaginosos"
$ws.Range("I876").Value = 100106002
$ws.Range("J876").Value = "Palta"
$ws.Range("K876").Value = "Hass"
$ws.Range("L876").Value = "Primera"
$ws.Range("M876").Value = 270
$ws.Range("N876").Value = 3000
$ws.Range("O876").Value = 3200
$ws.Range("P876").Value = 3111
$ws.Range("Q876").Value = "`$/kilo (en bandeja de 18 kilos)"
$ws.Range("R876").Value = "Provincia de Quillota"
$ws.Range("S876").Value = 3111
$ws.Range("T876").Value = 1

# Row 877 (new): Hass / Segunda
$ws.Range("A877").Value = 10
$ws.Range("B877").Value = "Vega Modelo de Temuco"
$ws.Range("C877").Value = "La Araucanía"
$ws.Range("D877").Value = 44585
$ws.Range("E877").Value = 9
$ws.Range("F877").Value = "Fruta"
$ws.Range("G877").Value = 100106
$ws.Range("H877").Value = "Oleaginosos"
$ws.Range("I877").Value = 100106002
$ws.Range("J877").Value = "Palta"
$ws.Range("K877").Value = "Hass"
$ws.Range("L877").Value = "Segunda"
$ws.Range("M877").Value = 100
$ws.Range("N877").Value = 2500
$ws.Range("O877").Value = 2500
$ws.Range("P877").Value = 2500
$ws.Range("Q877").Value = "`$/kilo (en bandeja de 18 kilos)"
$ws.Range("R877").Value = "Provincia de Quillota"
$ws.Range("S877").Value = 2500
$ws.Range("T877").Value = 1
